# Add a new worksheet "NewCarsTest" with car brand/title test data,
# and make it the active sheet (mirrors the authored diff).

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
# (tab order: LoginTest, NewCarsTest).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "NewCarsTest"

# Header row
$newSheet.Range("A1").Value = "carBrand"
$newSheet.Range("B1").Value = "carTitle"

# Data rows
$newSheet.Range("A2").Value = "BMW"
$newSheet.Range("B2").Value = "BMW Cars"

$newSheet.Range("A3").Value = "Hyundai"
$newSheet.Range("B3").Value = "Hyundai Cars"

$newSheet.Range("A4").Value = "Toyota"
$newSheet.Range("B4").Value = "Toyota Cars"

$newSheet.Range("A5").Value = "Honda"
$newSheet.Range("B5").Value = "Honda Cars"

# Make the new sheet the active/selected tab.
$newSheet.Select()
$newSheet.Range("I15").Select()
